$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format D cells whose new values are unambiguous numbers as Text,
# so Excel stores them as strings (matching source data) rather than
# auto-converting to floating point numbers.
$numericLookingCells = @('D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D19', 'D20', 'D22', 'D23', 'D27', 'D28', 'D29', 'D30', 'D31', 'D34', 'D37', 'D45', 'D46', 'D47', 'D48', 'D51')
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '43.203.81'
$ws.Range('E2').Value = '  +1.63%  '
# Row 3
$ws.Range('D3').Value = '2.304.89'
$ws.Range('E3').Value = '  +1.00%  '
# Row 4
$ws.Range('E4').Value = '  -0.03%  '
# Row 5
$ws.Range('D5').Value = '302.09'
$ws.Range('E5').Value = '  +0.76%  '
# Row 6
$ws.Range('D6').Value = '100.67'
$ws.Range('E6').Value = '  +5.48%  '
# Row 7
$ws.Range('D7').Value = '0.504'
$ws.Range('E7').Value = '  +0.84%  '
# Row 8
$ws.Range('E8').Value = '  +0.09%  '
# Row 9
$ws.Range('D9').Value = '0.516'
$ws.Range('E9').Value = '  +5.35%  '
# Row 10
$ws.Range('D10').Value = '36.59'
$ws.Range('E10').Value = '  +10.31%  '
# Row 11
$ws.Range('D11').Value = '0.0794'
$ws.Range('E11').Value = '  +0.73%  '
# Row 12
$ws.Range('E12').Value = '  +11.95%  '
# Row 13
$ws.Range('E13').Value = '  +1.61%  '
# Row 14
$ws.Range('E14').Value = '  +3.70%  '
# Row 15
$ws.Range('D15').Value = '2.666.43'
$ws.Range('E15').Value = '  +1.02%  '
# Row 16
$ws.Range('D16').Value = '2.309.01'
$ws.Range('E16').Value = '  +1.48%  '
# Row 17
$ws.Range('E17').Value = '  +1.70%  '
# Row 18
$ws.Range('D18').Value = '43.076.53'
$ws.Range('E18').Value = '  +1.51%  '
# Row 19
$ws.Range('D19').Value = '12.70'
$ws.Range('E19').Value = '  +11.33%  '
# Row 20
$ws.Range('D20').Value = '6.22'
$ws.Range('E20').Value = '  +4.40%  '
# Row 21
$ws.Range('D21').Value = '0.0₃0908'
$ws.Range('E21').Value = '  +1.76%  '
# Row 22
$ws.Range('D22').Value = '68.16'
$ws.Range('E22').Value = '  +2.22%  '
# Row 23
$ws.Range('D23').Value = '236.54'
$ws.Range('E23').Value = '  +0.57%  '
# Row 24
$ws.Range('E24').Value = '  +15.26%  '
# Row 25
$ws.Range('E25').Value = '  +0.33%  '
# Row 26
$ws.Range('E26').Value = '  +0.75%  '
# Row 27
$ws.Range('D27').Value = '25.14'
$ws.Range('E27').Value = '  +3.78%  '
# Row 28
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '2.29'
$ws.Range('E28').Value = '  +5.64%  '
# Row 29
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = '35.01'
$ws.Range('E29').Value = '  +4.42%  '
# Row 30
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '168.19'
$ws.Range('E30').Value = '  +0.92%  '
# Row 31
$ws.Range('D31').Value = '9.20'
$ws.Range('E31').Value = '  +1.63%  '
# Row 32
$ws.Range('E32').Value = '  -0.05%  '
# Row 33
$ws.Range('E33').Value = '  +2.52%  '
# Row 34
$ws.Range('D34').Value = '17.84'
$ws.Range('E34').Value = '  +6.00%  '
# Row 35
$ws.Range('E35').Value = '  -0.34%  '
# Row 36
$ws.Range('E36').Value = '  -0.47%  '
# Row 37
$ws.Range('D37').Value = '0.0698'
$ws.Range('E37').Value = '  +1.44%  '
# Row 38
$ws.Range('E38').Value = '  +1.92%  '
# Row 39
$ws.Range('E39').Value = '  +4.16%  '
# Row 40
$ws.Range('E40').Value = '  +2.03%  '
# Row 41
$ws.Range('E41').Value = '  +0.90%  '
# Row 43
$ws.Range('D43').Value = '1.990.30'
$ws.Range('E43').Value = '  +2.07%  '
# Row 44
$ws.Range('E44').Value = '  +4.61%  '
# Row 45
$ws.Range('D45').Value = '10.15'
$ws.Range('E45').Value = '  +4.99%  '
# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '17.73'
$ws.Range('E46').Value = '  +2.31%  '
# Row 47
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '2.92'
$ws.Range('E47').Value = '  +4.25%  '
# Row 48
$ws.Range('D48').Value = '55.85'
$ws.Range('E48').Value = '  +6.73%  '
# Row 49
$ws.Range('E49').Value = '  +5.94%  '
# Row 50
$ws.Range('D50').Value = '2.533.74'
$ws.Range('E50').Value = '  +0.91%  '
# Row 51
$ws.Range('D51').Value = '4.53'
$ws.Range('E51').Value = '  +0.24%  '

# Restore default style on the pre-formatted cells so no stray
# cell-level style survives beyond what the diff expects.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
